$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 3500333.8
$ws.Range("I12").Value = 5000500.5
$ws.Range("J12").Value = 500000
$ws.Range("K12").Value = 5000500.5
$ws.Range("L12").Value = 500000
$ws.Range("M12").Value = -5000330.5
$ws.Range("N12").Value = -500340
$ws.Range("H20").Value = 26708
$ws.Range("I20").Value = 5050
$ws.Range("J20").Value = 70024
$ws.Range("K20").Value = 5050
$ws.Range("L20").Value = 70024
$ws.Range("M20").Value = -4820
$ws.Range("N20").Value = -70484
$ws.Range("H32").Value = 1703.6
$ws.Range("J32").Value = 1150.2307
$ws.Range("L32").Value = 1150.2307
$ws.Range("N32").Value = -1802.2307
$ws.Range("H33").Value = 143018.86
$ws.Range("J33").Value = 164
$ws.Range("L33").Value = 164
$ws.Range("N33").Value = -622
$ws.Range("H35").Value = 26708
$ws.Range("I35").Value = 5050
$ws.Range("J35").Value = 70024
$ws.Range("K35").Value = 5050
$ws.Range("L35").Value = 70024
$ws.Range("M35").Value = -4671
$ws.Range("N35").Value = -70782
$ws.Range("H113").Value = 2495.64
$ws.Range("I113").Value = 2724.3125
$ws.Range("K113").Value = 2724.3125
$ws.Range("M113").Value = 529.6875
$ws.Range("H129").Value = 959.15625
$ws.Range("J129").Value = 1112.1154
$ws.Range("L129").Value = 3336.3462
$ws.Range("N129").Value = -13336.3462
$ws.Range("H132").Value = 314528.62
$ws.Range("I132").Value = 2140.1738
$ws.Range("K132").Value = 6420.5214
$ws.Range("M132").Value = -3890.5214
$ws.Range("H138").Value = 2764.41
$ws.Range("I138").Value = 872.8125
$ws.Range("J138").Value = 3124.7144
$ws.Range("K138").Value = 2618.4375
$ws.Range("L138").Value = 9374.143199999999
$ws.Range("M138").Value = 2521.5625
$ws.Range("N138").Value = -19654.1432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2323.4666
$ws.Range("I2").Value = 2313.4375
$ws.Range("J2").Value = 2348.1538
$ws.Range("K2").Value = 2313.4375
$ws.Range("L2").Value = 2348.1538
$ws.Range("M2").Value = -2200.4375
$ws.Range("N2").Value = -2574.1538
$ws.Range("H116").Value = 2323.4666
$ws.Range("I116").Value = 2313.4375
$ws.Range("J116").Value = 2348.1538
$ws.Range("K116").Value = 2313.4375
$ws.Range("L116").Value = 2348.1538
$ws.Range("M116").Value = -19.4375
$ws.Range("N116").Value = -6936.1538
$ws.Range("H132").Value = 2027.6666
$ws.Range("I132").Value = 1357.8182
$ws.Range("J132").Value = 2764.5
$ws.Range("K132").Value = 4073.4546
$ws.Range("L132").Value = 8293.5
$ws.Range("M132").Value = -1543.4546
$ws.Range("N132").Value = -13353.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2323.4666
$ws.Range("I3").Value = 2313.4375
$ws.Range("J3").Value = 2348.1538
$ws.Range("K3").Value = 2313.4375
$ws.Range("L3").Value = 2348.1538
$ws.Range("M3").Value = -2199.4375
$ws.Range("N3").Value = -2576.1538
$ws.Range("H135").Value = 72353
$ws.Range("J135").Value = 72353
$ws.Range("L135").Value = 72353
$ws.Range("N135").Value = -82493
$ws.Range("H137").Value = 53030.77
$ws.Range("J137").Value = 53030.77
$ws.Range("L137").Value = 53030.77
$ws.Range("N137").Value = -63230.77

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 289.29413
$ws.Range("I22").Value = 286.76923
$ws.Range("J22").Value = 297.5
$ws.Range("K22").Value = 286.76923
$ws.Range("L22").Value = 297.5
$ws.Range("M22").Value = 63.23077000000001
$ws.Range("N22").Value = -997.5
$ws.Range("H31").Value = 2923.6667
$ws.Range("I31").Value = 2969.85
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 2969.85
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -2674.85
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 2923.6667
$ws.Range("I34").Value = 2969.85
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 2969.85
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -2767.85
$ws.Range("N34").Value = -2404
$ws.Range("H99").Value = 402137.56
$ws.Range("I99").Value = 557135.6
$ws.Range("J99").Value = 3571.1428
$ws.Range("K99").Value = 557135.6
$ws.Range("L99").Value = 3571.1428
$ws.Range("M99").Value = -555637.6
$ws.Range("N99").Value = -6567.1428
$ws.Range("H126").Value = 402137.56
$ws.Range("I126").Value = 557135.6
$ws.Range("J126").Value = 3571.1428
$ws.Range("K126").Value = 1671406.8
$ws.Range("L126").Value = 10713.4284
$ws.Range("M126").Value = -1668936.8
$ws.Range("N126").Value = -15653.4284
$ws.Range("H132").Value = 2458.5
$ws.Range("I132").Value = 1902.15
$ws.Range("J132").Value = 3385.75
$ws.Range("K132").Value = 5706.450000000001
$ws.Range("L132").Value = 10157.25
$ws.Range("M132").Value = -3176.450000000001
$ws.Range("N132").Value = -15217.25
$ws.Range("H134").Value = 3695.1667
$ws.Range("I134").Value = 2993.3076
$ws.Range("J134").Value = 5520
$ws.Range("K134").Value = 8979.9228
$ws.Range("L134").Value = 16560
$ws.Range("M134").Value = -6444.9228
$ws.Range("N134").Value = -21630
$ws.Range("H140").Value = 53924
$ws.Range("J140").Value = 53924
$ws.Range("L140").Value = 53924
$ws.Range("N140").Value = -64284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1470.862
$ws.Range("I5").Value = 1165.8334
$ws.Range("K5").Value = 3497.5002
$ws.Range("M5").Value = -3385.5002
$ws.Range("H106").Value = 6598.3335
$ws.Range("I106").Value = 1100
$ws.Range("J106").Value = 7698
$ws.Range("K106").Value = 3300
$ws.Range("L106").Value = 23094
$ws.Range("M106").Value = -2354
$ws.Range("N106").Value = -24986
$ws.Range("H122").Value = 345659.34
$ws.Range("I122").Value = 544.94116
$ws.Range("J122").Value = 834571.4399999999
$ws.Range("K122").Value = 4904.47044
$ws.Range("L122").Value = 7511142.959999999
$ws.Range("M122").Value = -2454.47044
$ws.Range("N122").Value = -7516042.959999999
$ws.Range("H135").Value = 1470.862
$ws.Range("I135").Value = 1165.8334
$ws.Range("K135").Value = 10492.5006
$ws.Range("M135").Value = -7957.500599999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 150
$ws.Range("J4").Value = 150
$ws.Range("L4").Value = 150
$ws.Range("N4").Value = -374
$ws.Range("H5").Value = 8285
$ws.Range("J5").Value = 8285
$ws.Range("L5").Value = 8285
$ws.Range("N5").Value = -8509
$ws.Range("H70").Value = 4506.905
$ws.Range("I70").Value = 4049.9285
$ws.Range("K70").Value = 4049.9285
$ws.Range("M70").Value = -3779.9285
$ws.Range("H73").Value = 4506.905
$ws.Range("I73").Value = 4049.9285
$ws.Range("K73").Value = 4049.9285
$ws.Range("M73").Value = -3113.9285
$ws.Range("H113").Value = 1268.5
$ws.Range("I113").Value = 1268.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1268.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 901.5
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 22384.166
$ws.Range("J136").Value = 22384.166
$ws.Range("L136").Value = 67152.49800000001
$ws.Range("N136").Value = -72252.49800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1381.909
$ws.Range("I46").Value = 1650.1666
$ws.Range("J46").Value = 1060
$ws.Range("K46").Value = 1650.1666
$ws.Range("L46").Value = 1060
$ws.Range("M46").Value = -1462.1666
$ws.Range("N46").Value = -1436
$ws.Range("H55").Value = 494.14285
$ws.Range("J55").Value = 494.75
$ws.Range("L55").Value = 494.75
$ws.Range("N55").Value = -840.75
$ws.Range("H136").Value = 3005.6365
$ws.Range("I136").Value = 1306.2
$ws.Range("K136").Value = 3918.6
$ws.Range("M136").Value = -1368.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 10000
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("N8").Value = -10280
$ws.Range("H70").Value = 15000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 15000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 15000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 15000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 15000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 15000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -17184
